# Contactos.xlsx — form submission now works correctly; update the
# "Sr" -> "Hno" label and clean up the leftover underline style on D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Sr" role label to "Hno" everywhere it appears
#    (the "Genero" column for both contacts currently reads "Sr").
[void]$ws.Cells.Replace("Sr", "Hno", 1)

# 2. D8 still carries the old underlined placeholder style; make it match
#    D7's (yellow-fill, non-underlined) formatting by copying D7's format
#    onto D8, then keep the row height consistent with the rest of the sheet.
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 13.8

# 3. Nudge the neighbouring empty rows so the sheet's used area stretches
#    down to row 11, matching the rest of the form's row rhythm.
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(9).RowHeight = 13.8
$ws.Rows.Item(10).RowHeight = 13.8
$ws.Rows.Item(11).RowHeight = 13.8

# 4. Leave the cursor on C4, where the editor left off.
[void]$ws.Range("C4").Select()
